$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: Volume / Number and the reporting week date range
$ws.Range("A8").Value = "Volume 32   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/15/2025  Through  12/21/2025"


# --- Crime statistics table updates (rows 15-28, 31) ---
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("L15").Value = -14.285714285714
$ws.Range("C16").Value = 2
$ws.Range("D16").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -42.857142857142
$ws.Range("I16").Value = 160
$ws.Range("J16").Value = 201
$ws.Range("K16").Value = -20.398009950248
$ws.Range("L16").Value = -26.605504587156
$ws.Range("M16").Value = 13.475177304964
$ws.Range("N16").Value = -88.139362490733
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -25
$ws.Range("I17").Value = 205
$ws.Range("J17").Value = 198
$ws.Range("K17").Value = 3.535353535353
$ws.Range("L17").Value = 13.888888888888
$ws.Range("M17").Value = 120.430107526882
$ws.Range("N17").Value = -28.571428571428
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -36.363636363636
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = -54.838709677419
$ws.Range("I18").Value = 271
$ws.Range("J18").Value = 247
$ws.Range("K18").Value = 9.716599190283
$ws.Range("L18").Value = 8.835341365461
$ws.Range("M18").Value = 17.826086956521
$ws.Range("N18").Value = -90.816672314469
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = 42
$ws.Range("E19").Value = -23.809523809523
$ws.Range("F19").Value = 117
$ws.Range("G19").Value = 143
$ws.Range("H19").Value = -18.181818181818
$ws.Range("I19").Value = 1573
$ws.Range("J19").Value = 1638
$ws.Range("K19").Value = -3.968253968253
$ws.Range("L19").Value = -6.424747174301
$ws.Range("M19").Value = 26.243980738362
$ws.Range("N19").Value = -57.220560239325
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 60
$ws.Range("J20").Value = 86
$ws.Range("K20").Value = -30.232558139534
$ws.Range("L20").Value = -59.731543624161
$ws.Range("M20").Value = -31.818181818181
$ws.Range("N20").Value = -98.225377107364
$ws.Range("C21").Value = 47
$ws.Range("D21").Value = 65
$ws.Range("E21").Value = -27.692307692307
$ws.Range("F21").Value = 161
$ws.Range("G21").Value = 215
$ws.Range("H21").Value = -25.116279069767
$ws.Range("I21").Value = 2281
$ws.Range("J21").Value = 2388
$ws.Range("K21").Value = -4.480737018425
$ws.Range("L21").Value = -8.540497193263
$ws.Range("M21").Value = 25.744211686879
$ws.Range("N21").Value = -80.470890410958
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 29
$ws.Range("J22").Value = 43
$ws.Range("K22").Value = -32.558139534883
$ws.Range("L22").Value = -34.090909090909
$ws.Range("M22").Value = 11.538461538461
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 75
$ws.Range("I23").Value = 35
$ws.Range("J23").Value = 32
$ws.Range("K23").Value = 9.375
$ws.Range("L23").Value = 6.060606060606
$ws.Range("M23").Value = 45.833333333333
$ws.Range("C24").Value = 42
$ws.Range("D24").Value = 68
$ws.Range("E24").Value = -38.235294117647
$ws.Range("F24").Value = 171
$ws.Range("G24").Value = 283
$ws.Range("H24").Value = -39.575971731448
$ws.Range("I24").Value = 2723
$ws.Range("J24").Value = 3304
$ws.Range("K24").Value = -17.584745762711
$ws.Range("L24").Value = -9.263578807064
$ws.Range("M24").Value = 66.340867440439
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 59
$ws.Range("E25").Value = -62.711864406779
$ws.Range("F25").Value = 129
$ws.Range("G25").Value = 257
$ws.Range("H25").Value = -49.805447470817
$ws.Range("I25").Value = 2225
$ws.Range("J25").Value = 2898
$ws.Range("K25").Value = -23.222912353347
$ws.Range("L25").Value = -15.527714502657
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = -32.352941176470
$ws.Range("I26").Value = 350
$ws.Range("J26").Value = 367
$ws.Range("K26").Value = -4.632152588555
$ws.Range("L26").Value = 1.744186046511
$ws.Range("M26").Value = 4.477611940298
$ws.Range("C27").Value = "'0"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("L27").Value = -42.307692307692
$ws.Range("D28").Value = "'0"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "'***.*"
$ws.Range("E27").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 98
$ws.Range("K28").Value = -9.259259259259
$ws.Range("L28").Value = 13.953488372093
$ws.Range("L31").Value = -48.275862068965
